$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two product-availability labels (M column)
$ws.Range("M71").Value = "Reiswaffeln Fair Trade Milchschokolade - Online kein Bestand 1.30 Schweizer Franken"
$ws.Range("M95").Value = "Betty Bossi Naturaplan Bio Blätterteig ausgewallt - Online kein Bestand 2.40 Schweizer Franken"

# Update the scrape timestamp (O column) for every data row, 2 through 398
$newTimestamp = "2023-01-08 20:49:47"
for ($row = 2; $row -le 398; $row++) {
    $ws.Cells.Item($row, 15).Value = $newTimestamp
}
